# forests-scraped.xlsx update - 2025-12-02 12:20
# 1) The 7 rows currently on the "New" sheet (rows 2-8) are archived onto the
#    "Previously added" sheet (appended as rows 307-313), hyperlinks included.
# 2) The "New" sheet is cleared out (rows 5-8 removed) and rows 2-4 are
#    replaced with 3 freshly scraped listings (with their own hyperlinks).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# xlPasteValues / xlPasteFormats constants used with PasteSpecial below.
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Step 1: archive the current "New" rows 2-8 onto the end of "Previously added"
# ---------------------------------------------------------------------------

$firstNewRow = 307
$lastOldRow  = 8

# Old hyperlink targets living on "New" rows 2..8 (rId1..rId7), in order.
$oldLinkTargets = @(
    "https://www.ss.com/msg/lv/real-estate/wood/aluksne-and-reg/aluksne/gnedm.html",
    "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/pildas-pag/bbexhf.html",
    "https://www.ss.com/msg/lv/real-estate/wood/ogre-and-reg/taurupes-pag/cfmhh.html",
    "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/griskanu-pag/cghxgj.html",
    "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/kaunatas-pag/bgnoig.html",
    "https://www.ss.com/msg/lv/real-estate/wood/saldus-and-reg/saldus-pag/lklod.html",
    "https://www.ss.com/msg/lv/real-estate/wood/saldus-and-reg/saldus/jchdj.html"
)

for ($oldRow = 2; $oldRow -le $lastOldRow; $oldRow++) {
    $destRow = $firstNewRow + ($oldRow - 2)

    $srcRange  = $ws2.Range("A" + $oldRow + ":F" + $oldRow)
    $destRange = $ws1.Range("A" + $destRow + ":F" + $destRow)

    # Copy values first (keeps text cells that look numeric, e.g. cadastral
    # numbers, stored as text instead of being re-parsed as numbers), then
    # copy the formatting on top of it (keeps the existing style indices
    # instead of minting new ones).
    $srcRange.Copy()
    $destRange.PasteSpecial($xlPasteValues)
    $srcRange.Copy()
    $destRange.PasteSpecial($xlPasteFormats)

    $destCellA = $ws1.Range("A" + $destRow)
    $ws1.Hyperlinks.Add($destCellA, $oldLinkTargets[$oldRow - 2])

    # Hyperlinks.Add stamps its own built-in "Hyperlink" style onto the cell -
    # put the original column-A style back on top of it.
    $srcRange.Copy()
    $destRange.PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# Step 2: wipe the old hyperlinks + rows 5-8 from "New"
# ---------------------------------------------------------------------------

$ws2.Range("A2:A8").Hyperlinks.Delete()
$ws2.Rows("5:8").Delete()

# ---------------------------------------------------------------------------
# Step 3: write the 3 newly scraped listings onto "New" rows 2-4
# ---------------------------------------------------------------------------

$newRows = @(
    @{
        Link   = "https://www.ss.com/msg/lv/real-estate/wood/aluksne-and-reg/aluksne/gxkjp.html"
        Price  = "74 €"
        Region = "Alūksne un raj."
        Area   = "4 ha."
        Cad    = "36960010013"
        Date   = 45993.575694444444
    },
    @{
        Link   = "https://www.ss.com/msg/lv/real-estate/wood/ventspils-and-reg/zleku-pag/booep.html"
        Price  = "10 000 €"
        Region = "Ventspils un raj."
        Area   = "0.61 ha."
        Cad    = "98940030050"
        Date   = 45993.55416666667
    },
    @{
        Link   = "https://www.ss.com/msg/lv/real-estate/wood/ventspils-and-reg/zleku-pag/bbkhx.html"
        Price  = "10 000 €"
        Region = "Ventspils un raj."
        Area   = "3.98 ha."
        Cad    = "98940030021"
        Date   = 45993.55416666667
    }
)

# A reference range already carrying the correct per-column style (row 1's
# old row-2 data has long since moved, but the header-free data rows of
# "Previously added" still use the exact same style indices "New" needs).
$fmtRow = $ws1.Range("A307:F307")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row  = 2 + $i
    $data = $newRows[$i]

    $cellA = $ws2.Range("A" + $row)
    $cellB = $ws2.Range("B" + $row)
    $cellC = $ws2.Range("C" + $row)
    $cellD = $ws2.Range("D" + $row)
    $cellE = $ws2.Range("E" + $row)
    $cellF = $ws2.Range("F" + $row)

    $cellA.Value = $data.Link
    $cellC.Value = $data.Region
    $cellD.Value = $data.Area

    # "Price" and "Cad" look numeric (currency amount / cadastral id) - force
    # text formatting first so they are written as literal strings instead of
    # being auto-parsed into numbers.
    $cellB.NumberFormat = "@"
    $cellB.Value = $data.Price
    $cellE.NumberFormat = "@"
    $cellE.Value = $data.Cad

    $cellF.Value = $data.Date

    $ws2.Hyperlinks.Add($cellA, $data.Link)

    # Re-apply the canonical per-column styling (the NumberFormat tweaks
    # above, plus the freshly typed values, can otherwise leave the cells
    # with a mismatched style).
    $fmtRow.Copy()
    $ws2.Range("A" + $row + ":F" + $row).PasteSpecial($xlPasteFormats)
}
